$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 28: GenomeWeb link
$ws.Range("A28").Value = "https://www.genomeweb.com/molecular-diagnostics/qiagen-partnership-program-aims-build-out-clinical-dpcr-menu"
$ws.Range("B28").Value = "companion diagnostic"
$ws.Range("C28").Value = "Qiagen Partnership Program Aims to Build out Clinical dPCR Menu"

# New row 29: 360dx link
$ws.Range("A29").Value = "https://www.360dx.com/molecular-diagnostics/qiagen-partnership-program-aims-build-out-clinical-dpcr-menu"
$ws.Range("B29").Value = "companion diagnostic"
$ws.Range("C29").Value = "Qiagen Partnership Program Aims to Build out Clinical dPCR Menu"

# Add hyperlinks matching the URL text in column A, same as existing rows
$ws.Hyperlinks.Add($ws.Range("A28"), "https://www.genomeweb.com/molecular-diagnostics/qiagen-partnership-program-aims-build-out-clinical-dpcr-menu", "", "", "https://www.genomeweb.com/molecular-diagnostics/qiagen-partnership-program-aims-build-out-clinical-dpcr-menu")
$ws.Hyperlinks.Add($ws.Range("A29"), "https://www.360dx.com/molecular-diagnostics/qiagen-partnership-program-aims-build-out-clinical-dpcr-menu", "", "", "https://www.360dx.com/molecular-diagnostics/qiagen-partnership-program-aims-build-out-clinical-dpcr-menu")

$ws.Range("A28:A29").Style = "Hyperlink"
